# Generate Report for Handoff
# Updates the "b.md" file's status row (row 3) across the Overview, zh-cn
# and de-de sheets to reflect that a new handoff was generated.

$wb = $excel.ActiveWorkbook

function Set-HyperlinkDisplay($ws, $cellAddress, $newText) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $cellAddress) {
            $h.TextToDisplay = $newText
        }
    }
}

# ----- Overview sheet -----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-24-12 12:24:49"

# ----- zh-cn sheet -----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
Set-HyperlinkDisplay $zhcn '$D$3' "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-12 12:24:46"

# ----- de-de sheet -----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
Set-HyperlinkDisplay $dede '$D$3' "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-12 12:24:49"
